# Auto-generated edit script: applies scheduled market-price refresh
# to the Leve profit tables across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW.
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2221.5334  # H40
$ws.Cells.Item(40, 10).Value = 2550  # J40
$ws.Cells.Item(40, 12).Value = 2550  # L40
$ws.Cells.Item(40, 14).Value = -2900  # N40
$ws.Cells.Item(43, 8).Value = 675.5263  # H43
$ws.Cells.Item(43, 9).Value = 592.375  # I43
$ws.Cells.Item(43, 10).Value = 736  # J43
$ws.Cells.Item(43, 11).Value = 592.375  # K43
$ws.Cells.Item(43, 12).Value = 736  # L43
$ws.Cells.Item(43, 13).Value = -523.375  # M43
$ws.Cells.Item(43, 14).Value = -874  # N43
$ws.Cells.Item(53, 8).Value = 306.72223  # H53
$ws.Cells.Item(53, 9).Value = 204.27272  # I53
$ws.Cells.Item(53, 10).Value = 467.7143  # J53
$ws.Cells.Item(53, 11).Value = 204.27272  # K53
$ws.Cells.Item(53, 12).Value = 467.7143  # L53
$ws.Cells.Item(53, 13).Value = 432.72728  # M53
$ws.Cells.Item(53, 14).Value = -1741.7143  # N53
$ws.Cells.Item(61, 8).Value = 478.75  # H61
$ws.Cells.Item(61, 9).Value = 138.33333  # I61
$ws.Cells.Item(61, 10).Value = 1500  # J61
$ws.Cells.Item(61, 11).Value = 414.99999  # K61
$ws.Cells.Item(61, 12).Value = 4500  # L61
$ws.Cells.Item(61, 13).Value = -242.99999  # M61
$ws.Cells.Item(61, 14).Value = -4844  # N61
$ws.Cells.Item(76, 8).Value = 2806.4482  # H76
$ws.Cells.Item(76, 9).Value = 2523.762  # I76
$ws.Cells.Item(76, 10).Value = 3548.5  # J76
$ws.Cells.Item(76, 11).Value = 2523.762  # K76
$ws.Cells.Item(76, 12).Value = 3548.5  # L76
$ws.Cells.Item(76, 13).Value = -2208.762  # M76
$ws.Cells.Item(76, 14).Value = -4178.5  # N76
$ws.Cells.Item(79, 8).Value = 2806.4482  # H79
$ws.Cells.Item(79, 9).Value = 2523.762  # I79
$ws.Cells.Item(79, 10).Value = 3548.5  # J79
$ws.Cells.Item(79, 11).Value = 2523.762  # K79
$ws.Cells.Item(79, 12).Value = 3548.5  # L79
$ws.Cells.Item(79, 13).Value = -1431.762  # M79
$ws.Cells.Item(79, 14).Value = -5732.5  # N79
$ws.Cells.Item(125, 8).Value = 1866.6666  # H125
$ws.Cells.Item(125, 9).Value = 0  # I125
$ws.Cells.Item(125, 10).Value = 1866.6666  # J125
$ws.Cells.Item(125, 11).Value = 0  # K125
$ws.Cells.Item(125, 12).Value = 16799.9994  # L125
$ws.Cells.Item(125, 13).ClearContents()  # M125
$ws.Cells.Item(125, 14).Value = -21719.9994  # N125

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 1602.8572  # H88
$ws.Cells.Item(88, 9).Value = 1490  # I88
$ws.Cells.Item(88, 10).Value = 1687.5  # J88
$ws.Cells.Item(88, 11).Value = 1490  # K88
$ws.Cells.Item(88, 12).Value = 1687.5  # L88
$ws.Cells.Item(88, 13).Value = -1084  # M88
$ws.Cells.Item(88, 14).Value = -2499.5  # N88
$ws.Cells.Item(91, 8).Value = 1602.8572  # H91
$ws.Cells.Item(91, 9).Value = 1490  # I91
$ws.Cells.Item(91, 10).Value = 1687.5  # J91
$ws.Cells.Item(91, 11).Value = 1490  # K91
$ws.Cells.Item(91, 12).Value = 1687.5  # L91
$ws.Cells.Item(91, 13).Value = -86  # M91
$ws.Cells.Item(91, 14).Value = -4495.5  # N91

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(19, 8).Value = 20000  # H19
$ws.Cells.Item(19, 10).Value = 20000  # J19
$ws.Cells.Item(19, 12).Value = 20000  # L19
$ws.Cells.Item(19, 14).Value = -20346  # N19
$ws.Cells.Item(20, 8).Value = 2835.1177  # H20
$ws.Cells.Item(20, 9).Value = 3622  # I20
$ws.Cells.Item(20, 10).Value = 2593  # J20
$ws.Cells.Item(20, 11).Value = 3622  # K20
$ws.Cells.Item(20, 12).Value = 2593  # L20
$ws.Cells.Item(20, 13).Value = -3375  # M20
$ws.Cells.Item(20, 14).Value = -3087  # N20
$ws.Cells.Item(23, 8).Value = 5000  # H23
$ws.Cells.Item(23, 10).Value = 5000  # J23
$ws.Cells.Item(23, 12).Value = 5000  # L23
$ws.Cells.Item(23, 14).Value = -5566  # N23
$ws.Cells.Item(75, 8).Value = 13242.857  # H75
$ws.Cells.Item(75, 9).Value = 12540  # I75
$ws.Cells.Item(75, 10).Value = 15000  # J75
$ws.Cells.Item(75, 11).Value = 12540  # K75
$ws.Cells.Item(75, 12).Value = 15000  # L75
$ws.Cells.Item(75, 13).Value = -11604  # M75
$ws.Cells.Item(75, 14).Value = -16872  # N75
$ws.Cells.Item(78, 8).Value = 13242.857  # H78
$ws.Cells.Item(78, 9).Value = 12540  # I78
$ws.Cells.Item(78, 10).Value = 15000  # J78
$ws.Cells.Item(78, 11).Value = 37620  # K78
$ws.Cells.Item(78, 12).Value = 45000  # L78
$ws.Cells.Item(78, 13).Value = -32940  # M78
$ws.Cells.Item(78, 14).Value = -54360  # N78
$ws.Cells.Item(105, 8).Value = 2828.9333  # H105
$ws.Cells.Item(105, 9).Value = 1380  # I105
$ws.Cells.Item(105, 10).Value = 3051.8462  # J105
$ws.Cells.Item(105, 11).Value = 1380  # K105
$ws.Cells.Item(105, 12).Value = 3051.8462  # L105
$ws.Cells.Item(105, 13).Value = 367  # M105
$ws.Cells.Item(105, 14).Value = -6545.8462  # N105

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2854.1553  # H31
$ws.Cells.Item(31, 9).Value = 1854.4048  # I31
$ws.Cells.Item(31, 10).Value = 5478.5  # J31
$ws.Cells.Item(31, 11).Value = 1854.4048  # K31
$ws.Cells.Item(31, 12).Value = 5478.5  # L31
$ws.Cells.Item(31, 13).Value = -1559.4048  # M31
$ws.Cells.Item(31, 14).Value = -6068.5  # N31
$ws.Cells.Item(34, 8).Value = 2854.1553  # H34
$ws.Cells.Item(34, 9).Value = 1854.4048  # I34
$ws.Cells.Item(34, 10).Value = 5478.5  # J34
$ws.Cells.Item(34, 11).Value = 1854.4048  # K34
$ws.Cells.Item(34, 12).Value = 5478.5  # L34
$ws.Cells.Item(34, 13).Value = -1652.4048  # M34
$ws.Cells.Item(34, 14).Value = -5882.5  # N34
$ws.Cells.Item(59, 8).Value = 14169.0625  # H59
$ws.Cells.Item(59, 10).Value = 14169.0625  # J59
$ws.Cells.Item(59, 12).Value = 14169.0625  # L59
$ws.Cells.Item(59, 14).Value = -16459.0625  # N59
$ws.Cells.Item(86, 8).Value = 3169.1667  # H86
$ws.Cells.Item(86, 9).Value = 2835.6667  # I86
$ws.Cells.Item(86, 10).Value = 3502.6667  # J86
$ws.Cells.Item(86, 11).Value = 2835.6667  # K86
$ws.Cells.Item(86, 12).Value = 3502.6667  # L86
$ws.Cells.Item(86, 13).Value = -1712.6667  # M86
$ws.Cells.Item(86, 14).Value = -5748.6667  # N86
$ws.Cells.Item(89, 8).Value = 3169.1667  # H89
$ws.Cells.Item(89, 9).Value = 2835.6667  # I89
$ws.Cells.Item(89, 10).Value = 3502.6667  # J89
$ws.Cells.Item(89, 11).Value = 14178.3335  # K89
$ws.Cells.Item(89, 12).Value = 17513.3335  # L89
$ws.Cells.Item(89, 13).Value = -8562.333500000001  # M89
$ws.Cells.Item(89, 14).Value = -28745.3335  # N89

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 7.7  # H12
$ws.Cells.Item(12, 10).Value = 3.0625  # J12
$ws.Cells.Item(12, 12).Value = 9.1875  # L12
$ws.Cells.Item(12, 14).Value = -355.1875  # N12
$ws.Cells.Item(63, 8).Value = 4000  # H63
$ws.Cells.Item(63, 9).Value = 0  # I63
$ws.Cells.Item(63, 10).Value = 4000  # J63
$ws.Cells.Item(63, 11).Value = 0  # K63
$ws.Cells.Item(63, 12).Value = 12000  # L63
$ws.Cells.Item(63, 13).ClearContents()  # M63
$ws.Cells.Item(63, 14).Value = -13498  # N63
$ws.Cells.Item(66, 8).Value = 4000  # H66
$ws.Cells.Item(66, 9).Value = 0  # I66
$ws.Cells.Item(66, 10).Value = 4000  # J66
$ws.Cells.Item(66, 11).Value = 0  # K66
$ws.Cells.Item(66, 12).Value = 36000  # L66
$ws.Cells.Item(66, 13).ClearContents()  # M66
$ws.Cells.Item(66, 14).Value = -43488  # N66
$ws.Cells.Item(86, 8).Value = 192.6  # H86
$ws.Cells.Item(86, 9).Value = 190.75  # I86
$ws.Cells.Item(86, 10).Value = 200  # J86
$ws.Cells.Item(86, 11).Value = 572.25  # K86
$ws.Cells.Item(86, 12).Value = 600  # L86
$ws.Cells.Item(86, 13).Value = 613.75  # M86
$ws.Cells.Item(86, 14).Value = -2972  # N86
$ws.Cells.Item(89, 8).Value = 192.6  # H89
$ws.Cells.Item(89, 9).Value = 190.75  # I89
$ws.Cells.Item(89, 10).Value = 200  # J89
$ws.Cells.Item(89, 11).Value = 1716.75  # K89
$ws.Cells.Item(89, 12).Value = 1800  # L89
$ws.Cells.Item(89, 13).Value = 4211.25  # M89
$ws.Cells.Item(89, 14).Value = -13656  # N89

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 13483.5625  # H113
$ws.Cells.Item(113, 9).Value = 1542.5714  # I113
$ws.Cells.Item(113, 11).Value = 1542.5714  # K113
$ws.Cells.Item(113, 13).Value = 627.4286  # M113
$ws.Cells.Item(132, 8).Value = 3887.0571  # H132
$ws.Cells.Item(132, 9).Value = 4613.5293  # I132
$ws.Cells.Item(132, 10).Value = 3200.9443  # J132
$ws.Cells.Item(132, 11).Value = 13840.5879  # K132
$ws.Cells.Item(132, 12).Value = 9602.832900000001  # L132
$ws.Cells.Item(132, 13).Value = -11310.5879  # M132
$ws.Cells.Item(132, 14).Value = -14662.8329  # N132

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 7345.364  # H9
$ws.Cells.Item(9, 9).Value = 560  # I9
$ws.Cells.Item(9, 10).Value = 12999.833  # J9
$ws.Cells.Item(9, 11).Value = 560  # K9
$ws.Cells.Item(9, 12).Value = 12999.833  # L9
$ws.Cells.Item(9, 13).Value = -336  # M9
$ws.Cells.Item(9, 14).Value = -13447.833  # N9
$ws.Cells.Item(17, 8).Value = 10900  # H17
$ws.Cells.Item(17, 10).Value = 10900  # J17
$ws.Cells.Item(17, 12).Value = 10900  # L17
$ws.Cells.Item(17, 14).Value = -11240  # N17
$ws.Cells.Item(22, 8).Value = 617.4375  # H22
$ws.Cells.Item(22, 9).Value = 436.125  # I22
$ws.Cells.Item(22, 10).Value = 798.75  # J22
$ws.Cells.Item(22, 11).Value = 436.125  # K22
$ws.Cells.Item(22, 12).Value = 798.75  # L22
$ws.Cells.Item(22, 13).Value = -141.125  # M22
$ws.Cells.Item(22, 14).Value = -1388.75  # N22
$ws.Cells.Item(27, 8).Value = 617.4375  # H27
$ws.Cells.Item(27, 9).Value = 436.125  # I27
$ws.Cells.Item(27, 10).Value = 798.75  # J27
$ws.Cells.Item(27, 11).Value = 436.125  # K27
$ws.Cells.Item(27, 12).Value = 798.75  # L27
$ws.Cells.Item(27, 13).Value = -329.125  # M27
$ws.Cells.Item(27, 14).Value = -1012.75  # N27
$ws.Cells.Item(30, 8).Value = 1508  # H30
$ws.Cells.Item(30, 9).Value = 1508  # I30
$ws.Cells.Item(30, 11).Value = 1508  # K30
$ws.Cells.Item(30, 13).Value = -1400  # M30
$ws.Cells.Item(46, 8).Value = 1640.1  # H46
$ws.Cells.Item(46, 9).Value = 1000.5  # I46
$ws.Cells.Item(46, 10).Value = 1800  # J46
$ws.Cells.Item(46, 11).Value = 1000.5  # K46
$ws.Cells.Item(46, 12).Value = 1800  # L46
$ws.Cells.Item(46, 13).Value = -812.5  # M46
$ws.Cells.Item(46, 14).Value = -2176  # N46
$ws.Cells.Item(121, 8).Value = 24800  # H121
$ws.Cells.Item(121, 10).Value = 24800  # J121
$ws.Cells.Item(121, 12).Value = 24800  # L121
$ws.Cells.Item(121, 14).Value = -28294  # N121

